# Update the cryptos price/volume snapshot (GitHub Actions scrape refresh).
# For D29/D35/D48 the new price ends in a trailing zero (e.g. "7.10",
# "5.70", "0.0240") which Excel would otherwise normalize away if the
# text were auto-typed as a number (losing the trailing zero, e.g.
# "7.10" -> 7.1). A leading apostrophe forces those three cells to be
# entered as literal text so the digits match exactly, same as a user
# typing '7.10 into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.737.62"
$ws.Range("E2").Value = "  +2.48%  "

$ws.Range("D3").Value = "3.008.78"
$ws.Range("E3").Value = "  +2.33%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "510.89"
$ws.Range("E5").Value = "  +6.39%  "

$ws.Range("D6").Value = "139.24"
$ws.Range("E6").Value = "  +7.44%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  +4.86%  "

$ws.Range("E9").Value = "  +9.97%  "

$ws.Range("E10").Value = "  +9.71%  "

$ws.Range("D11").Value = "0.357"
$ws.Range("E11").Value = "  +3.72%  "

$ws.Range("E12").Value = "  +3.27%  "

$ws.Range("D13").Value = "3.522.19"
$ws.Range("E13").Value = "  +2.43%  "

$ws.Range("D14").Value = "25.81"
$ws.Range("E14").Value = "  +6.77%  "

$ws.Range("D15").Value = "0.0000158"
$ws.Range("E15").Value = "  +14.58%  "

$ws.Range("D16").Value = "56.784.81"
$ws.Range("E16").Value = "  +2.93%  "

$ws.Range("D17").Value = "3.007.80"
$ws.Range("E17").Value = "  +2.66%  "

$ws.Range("E18").Value = "  +7.93%  "

$ws.Range("D19").Value = "12.56"
$ws.Range("E19").Value = "  +6.29%  "

$ws.Range("D20").Value = "7.89"

$ws.Range("D21").Value = "331.65"
$ws.Range("E21").Value = "  +7.41%  "

$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("E23").Value = "  +6.43%  "

$ws.Range("D24").Value = "63.02"
$ws.Range("E24").Value = "  +5.96%  "

$ws.Range("D25").Value = "0.173"
$ws.Range("E25").Value = "  +11.19%  "

$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  +9.25%  "

$ws.Range("D28").Value = "6.74"
$ws.Range("E28").Value = "  +4.89%  "

$ws.Range("D29").Value = "'7.10"
$ws.Range("E29").Value = "  +10.65%  "

$ws.Range("E30").Value = "  +9.78%  "

$ws.Range("E31").Value = "  +8.26%  "

$ws.Range("D32").Value = "20.73"
$ws.Range("E32").Value = "  +8.58%  "

$ws.Range("D33").Value = "154.24"
$ws.Range("E33").Value = "  +5.01%  "

$ws.Range("E34").Value = "  +6.91%  "

$ws.Range("D35").Value = "'5.70"
$ws.Range("E35").Value = "  +2.59%  "

$ws.Range("D36").Value = "1.28"
$ws.Range("E36").Value = "  +1.98%  "

$ws.Range("D37").Value = "0.0681"
$ws.Range("E37").Value = "  +6.99%  "

$ws.Range("D38").Value = "24.18"
$ws.Range("E38").Value = "  +3.16%  "

$ws.Range("D39").Value = "3.040.94"
$ws.Range("E39").Value = "  +2.48%  "

$ws.Range("E40").Value = "  +3.46%  "

$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("E42").Value = "  +4.02%  "

$ws.Range("D43").Value = "2.277.09"
$ws.Range("E43").Value = "  +8.75%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "3.68"
$ws.Range("E44").Value = "  +5.84%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.42"
$ws.Range("E45").Value = "  +4.34%  "

$ws.Range("E46").Value = "  +1.61%  "

$ws.Range("D47").Value = "1.99"
$ws.Range("E47").Value = "  +21.06%  "

$ws.Range("D48").Value = "'0.0240"
$ws.Range("E48").Value = "  +6.01%  "

$ws.Range("D49").Value = "5.84"
$ws.Range("E49").Value = "  +6.25%  "

$ws.Range("D50").Value = "19.58"
$ws.Range("E50").Value = "  +5.21%  "

$ws.Range("D51").Value = "0.0875"
$ws.Range("E51").Value = "  +7.97%  "
